$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Nästa steg" notes (column D) for the rows that needed them ---
$ws.Range("D9").Value = "databasen uppdateras"
$ws.Range("D10").Value = "Databasen uppdateras när man lägger till ett bilmärke"
$ws.Range("D11").Value = "Databasen uppdateras när man lägger till ett land"
$ws.Range("D12").Value = "Man kan ta bort ett land från listan"
$ws.Range("D13").Value = "Landet försvinner från databasen om man raderar det"
$ws.Range("D14").Value = "Man får upp en msg box om man vill radera ett land"
$ws.Range("D15").Value = 'om man klickar "No" på msg boxen kommer man tillbaka'

# --- Test 10 turned out not to work ---
$ws.Range("B13").Value = "nej"

# --- Two brand-new test rows appended at the bottom ---
$ws.Range("A18").Value = "Test 16: Testa att man bara kan skriva in bokstäver på textrutor för namn osv."
$ws.Range("C18").Value = "textruta kan bara skriva in siffor"

$ws.Range("A19").Value = "Test 17: Testa att man bara kan skriva in siffror med en desimal på de textruror som behövs och att det kommer upp en msg box annars"

# --- Fill in "Fungerar?" (ja) answers for every other test, and the extra "E" column markers ---
$ws.Range("B6").Value = "ja"

$ws.Range("B7").Value = "ja"
$ws.Range("E7").Value = "ja"

$ws.Range("B8").Value = "ja"
$ws.Range("E8").Value = "ja"

$ws.Range("B9").Value = "ja"
$ws.Range("E9").Value = "ja"

$ws.Range("B10").Value = "ja"
$ws.Range("E10").Value = "ja"

$ws.Range("B11").Value = "ja"
$ws.Range("E11").Value = "ja"

$ws.Range("B12").Value = "ja"
$ws.Range("E12").Value = "ja"

$ws.Range("E13").Value = "ja"

$ws.Range("B14").Value = "ja"
$ws.Range("E14").Value = "ja"

$ws.Range("B15").Value = "ja"
$ws.Range("E15").Value = "ja"

$ws.Range("B16").Value = "ja"

$ws.Range("B17").Value = "ja"

# Match the borders/fill used by the rest of the table for the two new rows
$ws.Range("A18:A19").Borders.Item(7).LineStyle = 1
$ws.Range("A18:A19").Borders.Item(7).Weight = 4
$ws.Range("A18:A19").Borders.Item(10).LineStyle = 1
$ws.Range("A18:A19").Borders.Item(10).Weight = 2

$ws.Range("C18").Borders.Item(7).LineStyle = 1
$ws.Range("C18").Borders.Item(7).Weight = 2
$ws.Range("C18").Borders.Item(10).LineStyle = 1
$ws.Range("C18").Borders.Item(10).Weight = 4
$ws.Range("C18").HorizontalAlignment = -4108

# The column holding the test descriptions grew a bit wider after the new rows were added
$ws.Columns.Item(1).AutoFit()

# Selection was left on column B when the file was last saved
$ws.Columns.Item(2).Select()
